$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header date stamp - volatile NOW()
$ws.Range("N1").Formula = "=NOW()"

# Data rows 4-38: YEARFRAC age-in-years (col G) and renewal date H+365 (col I)
$ws.Range("G4:G38").Formula = "=YEARFRAC(F4,TODAY())"
$ws.Range("I4:I38").Formula = "=H4+365"

# Restore the view: selection on column I (matches the post-edit sheetView)
$ws.Range("I4:I38").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 7
} catch {}
